{"js": "const replacements = [\n  [\"694\u00d73=2082\", \"451\u00d72=902\"],\n  [\"102\u00d74=408\", \"464\u00d79=4176\"],\n  [\"577\u00d78=4616\", \"415\u00d78=3320\"],\n  [\"446\u00d74=1784\", \"625\u00d72=1250\"],\n  [\"276\u00d75=1380\", \"479\u00d77=3353\"],\n  [\"528\u00d79=4752\", \"344\u00d79=3096\"],\n  [\"435\u00d72=870\", \"443\u00d74=1772\"],\n  [\"799\u00d72=1598\", \"533\u00d78=4264\"],\n  [\"599\u00d72=1198\", \"693\u00d75=3465\"],\n  [\"691\u00d73=2073\", \"532\u00d76=3192\"],\n  [\"211\u00d72=422\", \"107\u00d75=535\"],\n  [\"380\u00d79=3420\", \"470\u00d79=4230\"],\n  [\"693\u00d73=2079\", \"656\u00d77=4592\"],\n  [\"115\u00d72=230\", \"264\u00d74=1056\"],\n  [\"350\u00d74=1400\", \"797\u00d73=2391\"],\n  [\"373\u00d72=746\", \"660\u00d77=4620\"],\n  [\"190\u00d73=570\", \"831\u00d78=6648\"],\n  [\"544\u00d79=4896\", \"193\u00d76=1158\"],\n  [\"673\u00d76=4038\", \"140\u00d75=700\"],\n  [\"885\u00d75=4425\", \"615\u00d78=4920\"],\n  [\"869\u00d72=1738\", \"802\u00d76=4812\"],\n  [\"463\u00d73=1389\", \"122\u00d75=610\"],\n  [\"803\u00d75=4015\", \"899\u00d73=2697\"],\n  [\"793\u00d78=6344\", \"240\u00d77=1680\"],\n  [\"776\u00d76=4656\", \"640\u00d72=1280\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n$replacements = @(\n    @{Old = \"694\u00d73=2082\"; New = \"451\u00d72=902\"},\n    @{Old = \"102\u00d74=408\"; New = \"464\u00d79=4176\"},\n    @{Old = \"577\u00d78=4616\"; New = \"415\u00d78=3320\"},\n    @{Old = \"446\u00d74=1784\"; New = \"625\u00d72=1250\"},\n    @{Old = \"276\u00d75=1380\"; New = \"479\u00d77=3353\"},\n    @{Old = \"528\u00d79=4752\"; New = \"344\u00d79=3096\"},\n    @{Old = \"435\u00d72=870\"; New = \"443\u00d74=1772\"},\n    @{Old = \"799\u00d72=1598\"; New = \"533\u00d78=4264\"},\n    @{Old = \"599\u00d72=1198\"; New = \"693\u00d75=3465\"},\n    @{Old = \"691\u00d73=2073\"; New = \"532\u00d76=3192\"},\n    @{Old = \"211\u00d72=422\"; New = \"107\u00d75=535\"},\n    @{Old = \"380\u00d79=3420\"; New = \"470\u00d79=4230\"},\n    @{Old = \"693\u00d73=2079\"; New = \"656\u00d77=4592\"},\n    @{Old = \"115\u00d72=230\"; New = \"264\u00d74=1056\"},\n    @{Old = \"350\u00d74=1400\"; New = \"797\u00d73=2391\"},\n    @{Old = \"373\u00d72=746\"; New = \"660\u00d77=4620\"},\n    @{Old = \"190\u00d73=570\"; New = \"831\u00d78=6648\"},\n    @{Old = \"544\u00d79=4896\"; New = \"193\u00d76=1158\"},\n    @{Old = \"673\u00d76=4038\"; New = \"140\u00d75=700\"},\n    @{Old = \"885\u00d75=4425\"; New = \"615\u00d78=4920\"},\n    @{Old = \"869\u00d72=1738\"; New = \"802\u00d76=4812\"},\n    @{Old = \"463\u00d73=1389\"; New = \"122\u00d75=610\"},\n    @{Old = \"803\u00d75=4015\"; New = \"899\u00d73=2697\"},\n    @{Old = \"793\u00d78=6344\"; New = \"240\u00d77=1680\"},\n    @{Old = \"776\u00d76=4656\"; New = \"640\u00d72=1280\"},\n)\n\nforeach ($item in $replacements) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Text = $item.Old\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Replacement.Text = $item.New\n    $range.Find.Execute($item.Old, $true, $false, $false, $false, $false, $true, 1, $false, $item.New, 2)\n}\n"}
